# First data upload all keyflows
# Re-sort the EURAL code table (A2:C26) in ascending EuralCode order. This
# moves the "04xxxx" rows above the "20xxxx" rows, drops the 200202 /
# "grond en stenen" row, and inserts a new 200123 row (Hazardous).
#
# Cell styles must stay pinned to their row position (the canonical diff
# never touches the `s="..."` attribute), so values are moved with
# PasteSpecial values-only, which carries the text/number type of the
# source cell without carrying its formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Stage the current A2:C26 values (values only, keeps t="s" text typing)
#    off to one side so source and destination ranges never collide while
#    we rewrite the table in place.
$ws.Range("A2:C26").Copy()
$ws.Range("E2").PasteSpecial(-4163)

# Map of target row -> staged source row (the row index inside E2:G26,
# i.e. same offset as the original A2:C26 row it came from).
# Row 16 (200123) has no source - it is a brand new row.
$srcRow = @{
  2  = 18
  3  = 19
  4  = 20
  5  = 21
  6  = 22
  7  = 23
  8  = 24
  9  = 25
  10 = 26
  11 = 2
  12 = 3
  13 = 4
  14 = 5
  15 = 6
  17 = 7
  18 = 8
  19 = 9
  20 = 10
  21 = 11
  22 = 12
  23 = 13
  24 = 15
  25 = 16
  26 = 17
}

foreach ($destRow in $srcRow.Keys) {
    $s = $srcRow[$destRow]
    $ws.Range("E$s`:G$s").Copy()
    $ws.Range("A$destRow").PasteSpecial(-4163)
}

# 2) Row 16 is the newly-introduced 200123 entry. Build the EuralCode text
#    value ("200123") via a helper cell whose format (only) is copied from
#    an existing Text-formatted EuralCode cell (A18, numFmtId 49) so typing
#    the digit string does not get coerced to the number 200123 - and so
#    doing this does not register a brand new cell style. Then
#    paste-special (values only) into A16 so A16 keeps its own
#    pre-existing style.
$helper = $ws.Range("Z1")
$ws.Range("A18").Copy()
$helper.PasteSpecial(-4122)
$helper.Value = "200123"
$helper.Copy()
$ws.Range("A16").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("B16").Value = "afgedankte apparatuur die chloorfluorkoolwaterstoffen bevat"
$ws.Range("C16").Value = "Hazardous"

# 3) Clean up the staging area so it doesn't show up as extra used range.
$ws.Range("E2:G26").Clear()

# 4) Restore the selected cell shown in the sheet view.
$ws.Range("C21").Select()
